$d = $word.ActiveDocument

# 1) "El sistema de" + (bookmark) + "be permitir inscribir a un socio en distintas
#    actividades..." was really one sentence split across a stray run/bookmark
#    boundary ("de" | "be" -> "debe"). Re-write it as a single run with the
#    correct merged text (this also removes the mis-placed _GoBack bookmark that
#    used to sit in the middle of "debe").
$d.Content.Find.Execute(
    "El sistema debe permitir inscribir a un socio en distintas actividades. El procedimiento será: dar de alta al socio, asociarle un abono, registrar las actividades que llevara a cabo.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El sistema debe permitir inscribir a un socio en distintas actividades. El procedimiento será: dar de alta al socio, asociarle un abono, registrar las actividades que llevara a cabo.",
    2) | Out-Null

# 2) Drop the stray leading space before "Notificaciones" further down the outline.
$d.Content.Find.Execute(
    " Notificaciones",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Notificaciones",
    2) | Out-Null

# 3) Put the _GoBack bookmark back at the very end of the document (the last,
#    empty paragraph right before the section break), which is where Word leaves
#    it after the edits above.
if (-not $d.Bookmarks.Exists("_GoBack")) {
    $lastPara = $d.Paragraphs.Last
    $d.Bookmarks.Add("_GoBack", $lastPara.Range) | Out-Null
}
